$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 10: header labels (written first so the shared-string table fills in
# the same order as the target: [0:100], [1:101], [0:10], [1:11], then the
# python command string)
$ws.Range("A10").Value = "MNIST, test indices [0:100], test time only"
$ws.Range("B10").Value = "MNIST, test indices[1:101], test time only"
$ws.Range("C10").Value = "MNIST, test indices [0:10], test time only"
$ws.Range("D10").Value = "MNIST, test indices [1:11], test time only"

# Row 9: the python command, styled like inline code (Consolas 8pt, GitHub-ish
# dark gray, left/vcenter aligned with a 1-level indent)
$ws.Range("A9").Value = "python mnist.py --device=cpu -n=1 --lr=.25 --sigma=1.3 -c=1.5 -b=250"
$ws.Range("A9").Font.Name = "Consolas"
$ws.Range("A9").Font.Family = 3
$ws.Range("A9").Font.Size = 8
$ws.Range("A9").Font.Color = 3025188
$ws.Range("A9").HorizontalAlignment = -4131
$ws.Range("A9").VerticalAlignment = -4108
$ws.Range("A9").IndentLevel = 1

# Rows 11-17: the new (smaller) test-time-only timing data
$data = @(
    @(1.8444, 1.8488, 1.8041, 1.3995),
    @(1.4799, 1.4223, 1.4115, 1.3873),
    @(1.8635, 1.4136, 1.3997, 1.3904),
    @(1.4968, 1.8388, 1.3881, 1.8224),
    @(1.5632, 1.8508, 1.8042, 1.8183),
    @(1.409, 1.8271, 1.8561, 1.8111),
    @(1.422, 1.8409, 1.3793, 1.7976)
)

$row = 11
foreach ($r in $data) {
    $ws.Cells.Item($row, 1).Value = $r[0]
    $ws.Cells.Item($row, 2).Value = $r[1]
    $ws.Cells.Item($row, 3).Value = $r[2]
    $ws.Cells.Item($row, 4).Value = $r[3]
    $row++
}

# Page setup: portrait orientation
$ws.PageSetup.Orientation = 1

# Final selection, matching the saved cursor position
$ws.Range("B18").Select()
